$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = "Riccardo Zeni"
$ws.Range("B50").Value = "Samuele Kettamier | SBARX"
$ws.Range("C50").Value = "Eduardo  Grazioli  | FC Savignano"
$ws.Range("D50").Value = "Mattia Baldessarini | Shark Attack"
$ws.Range("E50").Value = "Nadir Chtioui | MAI UNA GIOIA"
$ws.Range("F50").Value = "Roberto Barozzi | Demobusters"
